# This script permutes the per-row data (columns D, L, M, N, O, P, Q, R, S, T)
# across rows 2-18 of the active worksheet, per the target diff. Columns
# A, B, C, E, F, G, H, I, J, K are identical across all rows and are left
# untouched. The mapping below gives, for each destination row, the row
# whose current data should be copied into it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destinationRow -> sourceRow (based on the current/original row layout)
$mapping = @{
    2  = 10
    3  = 11
    4  = 12
    5  = 13
    6  = 16
    7  = 14
    8  = 15
    9  = 6
    10 = 7
    11 = 2
    12 = 3
    13 = 8
    14 = 9
    15 = 17
    16 = 18
    17 = 4
    18 = 5
}

$cols = @("D", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# Snapshot the current values for the columns that will move, for every
# source row, before any writes happen (writes must not affect reads).
$snapshot = @{}
foreach ($row in 2..18) {
    $rowData = @{}
    foreach ($col in $cols) {
        $rowData[$col] = $ws.Range("$col$row").Value()
    }
    $snapshot[$row] = $rowData
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $rowData = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $rowData[$col]
    }
}
